# Generate Report for Handoff
# Updates the handoff-priority flag ("ht") and the corresponding handoff
# timestamp columns for the rows that were just (re)generated for handoff,
# on the Overview, zh-cn, and de-de worksheets.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 12, 13, 14)

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-09-07 08:29:59"
}

# --- zh-cn sheet: "Priority" column (E) and "Latest Handoff Datetime" column (H) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-09-07 08:29:53"
}

# --- de-de sheet: "Priority" column (E) and "Latest Handoff Datetime" column (H) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-09-07 08:29:59"
}
